$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestSteps")
$ws2 = $wb.Worksheets.Item("TestData")

# --- TestSteps sheet: insert a new row 8 (shifts old rows 8,9 to 9,10), and
# append a new row 11 at the bottom ---
$ws1.Rows("8:8").Insert()

# Copy formatting (styles/borders) from row 7 into the new row 8 and the
# brand-new row 11 so they match the sheet's existing look.
$ws1.Range("A7:F7").Copy()
$ws1.Range("A8:F8").PasteSpecial(-4122)
$ws1.Range("A11:F11").PasteSpecial(-4122)

# Fill in row 11 first, then row 8, to mirror the original authoring order
# (keeps shared-string table ordering identical to the source edit).
$ws1.Range("A11").Value = "checkAccessibility"
$ws1.Range("B11").Value = "ListView_AfterDelete"
$ws1.Range("B8").Value = "DeletePopUp"

# --- TestData sheet: fill in the two NACUBO GL Account rows ---
$ws2.Range("B12").Value = "Enter Text to delete NACUBO GL Account Class Profile"
$ws2.Range("B11").Value = "Enter Text to delete NACUBO GL Account Category Profile"

$ws1.Range("A8").Value = "checkAccessibility"

# Match style of C9:C10 (s="6") for the new C11/C12 cells.
$ws2.Range("C9:C10").Copy()
$ws2.Range("C11:C12").PasteSpecial(-4122)

$ws2.Range("A11").Value = 10
$ws2.Range("C11").Value = 99
$ws2.Range("A12").Value = 11
$ws2.Range("C12").Value = "Adarsh"

# --- restore selections to match the target workbook state ---
$ws1.Range("B8").Select() | Out-Null
$ws2.Range("B11").Select() | Out-Null
